$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date cell I1
$ws.Range("I1").Value = "25/03/2023"

# Row 2
$ws.Range("C2").Value = 446
$ws.Range("D2").Value = 452
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 892
$ws.Range("J2").Value = -49.32735426008968

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 261
$ws.Range("D3").Value = 266
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 292
$ws.Range("J3").Value = -8.904109589041099

# Row 4
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 9
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 28.57142857142858

# Row 5
$ws.Range("C5").Value = 166
$ws.Range("D5").Value = 181
$ws.Range("F5").Value = 0
$ws.Range("I5").Value = 170
$ws.Range("J5").Value = 6.470588235294117

# Row 6
$ws.Range("C6").Value = 26
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = 39
$ws.Range("J6").Value = -33.33333333333334

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 49
$ws.Range("D7").Value = 51
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 68
$ws.Range("J7").Value = -25

# Row 8
$ws.Range("C8").Value = 159
$ws.Range("D8").Value = 160
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 19
$ws.Range("J8").Value = 742.1052631578947

# Row 9
$ws.Range("C9").Value = 23
$ws.Range("D9").Value = 23
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("I9").Value = 44
$ws.Range("J9").Value = -47.72727272727273

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 203
$ws.Range("D10").Value = 252
$ws.Range("E10").Value = 55
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7
$ws.Range("I10").Value = 394
$ws.Range("J10").Value = -36.04060913705583

# Row 11
$ws.Range("C11").Value = 161
$ws.Range("D11").Value = 163
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 218
$ws.Range("J11").Value = -25.22935779816514

# Row 12
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 341
$ws.Range("D12").Value = 453
$ws.Range("E12").Value = 25
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 84
$ws.Range("I12").Value = 774
$ws.Range("J12").Value = -41.47286821705426

# Row 13
$ws.Range("C13").Value = 62
$ws.Range("D13").Value = 69
$ws.Range("E13").Value = 6
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 486
$ws.Range("J13").Value = -85.80246913580247

# Row 14
$ws.Range("C14").Value = 511
$ws.Range("D14").Value = 584
$ws.Range("E14").Value = 30
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = 41
$ws.Range("I14").Value = 603
$ws.Range("J14").Value = -3.150912106135983

# Row 15
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 171
$ws.Range("D15").Value = 181
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 167
$ws.Range("J15").Value = 8.383233532934131

# Row 16
$ws.Range("C16").Value = 63
$ws.Range("D16").Value = 92
$ws.Range("E16").Value = 29
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = -23.33333333333333

# Row 17
$ws.Range("C17").Value = 39
$ws.Range("D17").Value = 37
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = -49.31506849315068

# Row 18
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = -66.66666666666667

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 9
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 0

# Row 20
$ws.Range("C20").Value = 28
$ws.Range("D20").Value = 28
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 1
$ws.Range("I20").Value = 43
$ws.Range("J20").Value = -34.88372093023256

